$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: Volume number 15 -> 16
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "16"

# Header: report week date range 4/7/2025-4/13/2025 -> 4/14/2025-4/20/2025
$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "4/14/2025"
$c9.Characters(47, 9).Text = "4/20/2025"

# --- Cells switching from numeric to text "0" (style copied from D14) ---
$ws.Range("D28").Value = "'0"
$ws.Range("D29").Value = "'0"
$ws.Range("D30").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D28,D29,D30").PasteSpecial(-4122)

# --- Cells switching from numeric to text "***.*" (style copied from E14) ---
$ws.Range("E28").Value = "***.*"
$ws.Range("E29").Value = "***.*"
$ws.Range("E30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E28,E29,E30").PasteSpecial(-4122)

# --- Cells switching from text "0" to numeric (style copied from D16) ---
$ws.Range("C15").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D16").Copy()
$ws.Range("C15,C27").PasteSpecial(-4122)

# --- Plain numeric value updates ---
$ws.Range("G14").Value = 3
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 133.333333333333
$ws.Range("L15").Value = 75
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -15.384615384615
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = -22.222222222222
$ws.Range("L16").Value = -24.324324324324
$ws.Range("M16").Value = -16.417910447761
$ws.Range("N16").Value = -75
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -13.793103448275
$ws.Range("I17").Value = 84
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = -38.686131386861
$ws.Range("L17").Value = -19.230769230769
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = -50.877192982456
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 27
$ws.Range("K18").Value = 40.740740740740
$ws.Range("L18").Value = 11.764705882352
$ws.Range("M18").Value = 26.666666666666
$ws.Range("N18").Value = -72.857142857142
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -53.846153846153
$ws.Range("F19").Value = 32
$ws.Range("H19").Value = -37.254901960784
$ws.Range("I19").Value = 112
$ws.Range("J19").Value = 152
$ws.Range("K19").Value = -26.315789473684
$ws.Range("L19").Value = -3.448275862068
$ws.Range("M19").Value = 93.103448275862
$ws.Range("N19").Value = 36.585365853658
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -47.058823529411
$ws.Range("I20").Value = 22
$ws.Range("J20").Value = 41
$ws.Range("K20").Value = -46.341463414634
$ws.Range("L20").Value = 29.411764705882
$ws.Range("M20").Value = 83.333333333333
$ws.Range("N20").Value = -74.712643678160
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -16.216216216216
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -25.619834710743
$ws.Range("I21").Value = 319
$ws.Range("J21").Value = 437
$ws.Range("K21").Value = -27.002288329519
$ws.Range("L21").Value = -9.116809116809
$ws.Range("M21").Value = 29.149797570850
$ws.Range("N21").Value = -55.571030640668
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 66.666666666666
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -7.692307692307
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 56
$ws.Range("K23").Value = -16.071428571428
$ws.Range("L23").Value = -2.083333333333
$ws.Range("M23").Value = 95.833333333333
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = -19.718309859154
$ws.Range("I24").Value = 229
$ws.Range("J24").Value = 272
$ws.Range("K24").Value = -15.808823529411
$ws.Range("L24").Value = -40.673575129533
$ws.Range("M24").Value = -21.575342465753
$ws.Range("C25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = -25
$ws.Range("I25").Value = 65
$ws.Range("J25").Value = 73
$ws.Range("K25").Value = -10.958904109589
$ws.Range("L25").Value = -65.240641711229
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -31.25
$ws.Range("F26").Value = 53
$ws.Range("G26").Value = 68
$ws.Range("H26").Value = -22.058823529411
$ws.Range("I26").Value = 184
$ws.Range("J26").Value = 217
$ws.Range("K26").Value = -15.207373271889
$ws.Range("L26").Value = 22.666666666666
$ws.Range("M26").Value = 49.593495934959
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 60
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("I28").Value = 17
$ws.Range("K28").Value = 13.333333333333
$ws.Range("L28").Value = -5.555555555555
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("L29").Value = -76.923076923076
$ws.Range("M29").Value = -62.5
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = -50
$ws.Range("L30").Value = -75
$ws.Range("M30").Value = -57.142857142857
